$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to text format so numeric-looking strings
# (e.g. "10.20", "0.4700") keep their exact original formatting
# instead of being auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.890.26"
$ws.Range("E2").Value = "  +0.53%  "

$ws.Range("D3").Value = "1.877.88"
$ws.Range("E3").Value = "  -0.04%  "

$ws.Range("E4").Value = "  +1.45%  "

$ws.Range("D5").Value = "334.62"
$ws.Range("E5").Value = "  +0.86%  "

$ws.Range("D6").Value = "1.015"
$ws.Range("E6").Value = "  +1.23%  "

$ws.Range("D7").Value = "0.4700"
$ws.Range("E7").Value = "  -0.38%  "

$ws.Range("D8").Value = "0.3912"
$ws.Range("E8").Value = "  -1.25%  "

$ws.Range("D9").Value = "46.72"
$ws.Range("E9").Value = "  -2.55%  "

$ws.Range("D10").Value = "0.07945"
$ws.Range("E10").Value = "  -1.13%  "

$ws.Range("D11").Value = "1.006"
$ws.Range("E11").Value = "  -1.75%  "

$ws.Range("D12").Value = "21.61"
$ws.Range("E12").Value = "  -1.06%  "

$ws.Range("D13").Value = "1.891.42"
$ws.Range("E13").Value = "  +0.24%  "

$ws.Range("D14").Value = "5.940"
$ws.Range("E14").Value = "  -0.45%  "

$ws.Range("D15").Value = "7.095"
$ws.Range("E15").Value = "  -0.83%  "

$ws.Range("D16").Value = "1.019"
$ws.Range("E16").Value = "  +1.45%  "

$ws.Range("D17").Value = "0.06782"
$ws.Range("E17").Value = "  +2.49%  "

$ws.Range("D18").Value = "87.38"
$ws.Range("E18").Value = "  +0.25%  "

$ws.Range("D19").Value = "0.00001045"
$ws.Range("E19").Value = "  -0.36%  "

$ws.Range("D20").Value = "17.02"
$ws.Range("E20").Value = "  -1.65%  "

$ws.Range("D21").Value = "1.015"
$ws.Range("E21").Value = "  +1.27%  "

$ws.Range("D22").Value = "27.912.22"
$ws.Range("E22").Value = "  +0.53%  "

$ws.Range("D23").Value = "5.473"
$ws.Range("E23").Value = "  -0.72%  "

$ws.Range("D24").Value = "10.95"
$ws.Range("E24").Value = "  -0.76%  "

$ws.Range("D25").Value = "2.354"
$ws.Range("E25").Value = "  +2.47%  "

$ws.Range("D26").Value = "2.097.30"
$ws.Range("E26").Value = "  -0.59%  "

$ws.Range("D27").Value = "159.94"
$ws.Range("E27").Value = "  +2.10%  "

$ws.Range("D28").Value = "19.89"
$ws.Range("E28").Value = "  -1.72%  "

$ws.Range("E29").Value = "  -0.58%  "

$ws.Range("D30").Value = "5.456"
$ws.Range("E30").Value = "  -2.49%  "

$ws.Range("D31").Value = "120.84"
$ws.Range("E31").Value = "  -1.36%  "

$ws.Range("D32").Value = "0.09530"
$ws.Range("E32").Value = "  -0.22%  "

$ws.Range("D33").Value = "0.9591"
$ws.Range("E33").Value = "  -1.06%  "

$ws.Range("D34").Value = "3.652"
$ws.Range("E34").Value = "  +0.69%  "

$ws.Range("D35").Value = "5.320"
$ws.Range("E35").Value = "  +0.34%  "

$ws.Range("D36").Value = "1.345"
$ws.Range("E36").Value = "  -7.57%  "

$ws.Range("D37").Value = "0.06103"
$ws.Range("E37").Value = "  -0.14%  "

$ws.Range("E38").Value = "  -1.01%  "

$ws.Range("D39").Value = "1.206"
$ws.Range("E39").Value = "  -2.09%  "

$ws.Range("D40").Value = "1.015"
$ws.Range("E40").Value = "  +1.28%  "

$ws.Range("D41").Value = "8.088"
$ws.Range("E41").Value = "  -1.21%  "

$ws.Range("D42").Value = "0.5918"
$ws.Range("E42").Value = "  -1.39%  "

$ws.Range("D43").Value = "0.1888"
$ws.Range("E43").Value = "  -0.99%  "

$ws.Range("D44").Value = "10.20"
$ws.Range("E44").Value = "  -0.50%  "

$ws.Range("D45").Value = "1.267"
$ws.Range("E45").Value = "  +1.20%  "

$ws.Range("D46").Value = "0.5620"
$ws.Range("E46").Value = "  -1.14%  "

$ws.Range("D47").Value = "12.11"
$ws.Range("E47").Value = "  -0.83%  "

$ws.Range("D48").Value = "3.403"
$ws.Range("E48").Value = "  +0.01%  "

$ws.Range("D49").Value = "1.917"
$ws.Range("E49").Value = "  -0.83%  "

$ws.Range("D50").Value = "0.06858"
$ws.Range("E50").Value = "  +0.50%  "

$ws.Range("D51").Value = "113.52"
$ws.Range("E51").Value = "  +1.67%  "
